$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 190, shifting existing rows 190-264 down to 191-265
$ws.Rows.Item(190).Insert()

# Populate the newly inserted row 190 with data
$ws.Cells.Item(190, 1).Value = "Marylebone - Monday Class (King Solomon Academy)"
$ws.Cells.Item(190, 2).Value = "October"
$ws.Cells.Item(190, 3).Value = 7
$ws.Cells.Item(190, 4).Value = "<3 months"
$ws.Cells.Item(190, 5).Value = 4
$ws.Cells.Item(190, 6).Value = 4
$ws.Cells.Item(190, 7).Value = 4
$ws.Cells.Item(190, 8).Value = 3
$ws.Cells.Item(190, 9).Value = 3
$ws.Cells.Item(190, 10).Value = 6
$ws.Cells.Item(190, 11).Value = 3
$ws.Cells.Item(190, 12).Value = 3
$ws.Cells.Item(190, 13).Value = 3
$ws.Cells.Item(190, 14).Value = 3

# Update window view / selection to mirror final saved state
$ws.Range("P195").Select()
$excel.ActiveWindow.ScrollRow = 177
$excel.ActiveWindow.ScrollColumn = 1
